$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price cells (column D) stay as text, matching the original inline-string formatting
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.266.05'
$ws.Range("E2").Value = '  -4.21%  '

$ws.Range("D3").Value = '1.855.39'
$ws.Range("E3").Value = '  -5.36%  '

$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -1.26%  '

$ws.Range("E5").Value = '  -0.36%  '

$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -1.00%  '

$ws.Range("D7").Value = '0.4501'
$ws.Range("E7").Value = '  -5.46%  '

$ws.Range("D8").Value = '0.3848'
$ws.Range("E8").Value = '  -4.77%  '

$ws.Range("D9").Value = '47.82'
$ws.Range("E9").Value = '  -11.46%  '

$ws.Range("D10").Value = '0.07875'
$ws.Range("E10").Value = '  -6.73%  '

$ws.Range("D11").Value = '1.016'
$ws.Range("E11").Value = '  -3.85%  '

$ws.Range("D12").Value = '21.34'
$ws.Range("E12").Value = '  -4.34%  '

$ws.Range("D13").Value = '1.856.52'
$ws.Range("E13").Value = '  -6.99%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '5.865'
$ws.Range("E14").Value = '  -4.83%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '7.142'
$ws.Range("E15").Value = '  -5.70%  '

$ws.Range("E16").Value = '  -1.33%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.00001030'
$ws.Range("E17").Value = '  -3.82%  '

$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").Value = '85.65'
$ws.Range("E18").Value = '  -5.59%  '

$ws.Range("D19").Value = '0.06519'
$ws.Range("E19").Value = '  -1.74%  '

$ws.Range("D20").Value = '16.93'
$ws.Range("E20").Value = '  -8.55%  '

$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  -1.07%  '

$ws.Range("D22").Value = '5.494'
$ws.Range("E22").Value = '  -6.11%  '

$ws.Range("D23").Value = '27.275.83'
$ws.Range("E23").Value = '  -4.37%  '

$ws.Range("D24").Value = '10.77'
$ws.Range("E24").Value = '  -5.80%  '

$ws.Range("E25").Value = '  -1.74%  '

$ws.Range("D26").Value = '2.070.60'
$ws.Range("E26").Value = '  -7.13%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '151.82'
$ws.Range("E27").Value = '  -2.44%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '19.74'
$ws.Range("E28").Value = '  -2.55%  '

$ws.Range("D29").Value = '2.060'
$ws.Range("E29").Value = '  -4.56%  '

$ws.Range("D30").Value = '5.483'
$ws.Range("E30").Value = '  -6.79%  '

$ws.Range("D31").Value = '120.21'
$ws.Range("E31").Value = '  -3.40%  '

$ws.Range("D32").Value = '0.9369'
$ws.Range("E32").Value = '  -4.31%  '

$ws.Range("D33").Value = '0.09270'
$ws.Range("E33").Value = '  -3.88%  '

$ws.Range("D34").Value = '1.465'
$ws.Range("E34").Value = '  +1.32%  '

$ws.Range("D35").Value = '3.566'
$ws.Range("E35").Value = '  -3.53%  '

$ws.Range("D36").Value = '5.293'
$ws.Range("E36").Value = '  -5.59%  '

$ws.Range("D37").Value = '0.02219'
$ws.Range("E37").Value = '  -4.67%  '

$ws.Range("D38").Value = '0.05976'
$ws.Range("E38").Value = '  -3.86%  '

$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = '1.204'
$ws.Range("E39").Value = '  -3.82%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '8.305'
$ws.Range("E40").Value = '  -9.23%  '

$ws.Range("D41").Value = '1.000'
$ws.Range("E41").Value = '  -1.07%  '

$ws.Range("D42").Value = '0.5903'
$ws.Range("E42").Value = '  -4.77%  '

$ws.Range("D43").Value = '0.1881'
$ws.Range("E43").Value = '  -1.46%  '

$ws.Range("D44").Value = '10.10'
$ws.Range("E44").Value = '  -9.29%  '

$ws.Range("D45").Value = '1.261'
$ws.Range("E45").Value = '  -7.03%  '

$ws.Range("D46").Value = '0.5616'
$ws.Range("E46").Value = '  -5.36%  '

$ws.Range("D47").Value = '11.85'
$ws.Range("E47").Value = '  -8.68%  '

$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").Value = '3.350'
$ws.Range("E48").Value = '  -1.69%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '1.916'
$ws.Range("E49").Value = '  -6.48%  '

$ws.Range("D50").Value = '0.06809'
$ws.Range("E50").Value = '  -0.08%  '

$ws.Range("D51").Value = '108.19'
$ws.Range("E51").Value = '  -2.62%  '
